$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B='Mistral-Small-3.2-24b'; C=2.659567901234568; D=1.802932098765432; E=2.709259259259259},
    @{Row=3; B='Claude-4.5-Sonnet'; C=2.629166666666667; D=2.614197530864198; E=2.918981481481481},
    @{Row=4; B='Claude-4.5-Haiku'; C=2.604166666666667; D=2.538271604938272; E=2.876543209876543},
    @{Row=5; B='Mistral-Small-24b-2501'; C=2.578240740740741; D=1.945833333333333; E=2.778703703703704},
    @{Row=6; B='Gemma-2-9b-It'; C=2.520061728395062; D=2.432098765432099; E=2.835802469135802},
    @{Row=7; B='Gemma-2-27b-It'; C=2.502006172839506; D=2.459567901234568; E=2.85787037037037},
    @{Row=8; B='Llama-3.3-8b'; C=2.459722222222222; D=1.977006172839507; E=2.678549382716049},
    @{Row=9; B='Claude-Sonnet-4'; C=2.440277777777778; D=2.316358024691358; E=2.830246913580247},
    @{Row=10; B='Claude-3.5-Sonnet'; C=2.394907407407407; D=2.624382716049383; E=2.953703703703704},
    @{Row=11; B='Claude-Opus-4'; C=2.364814814814815; D=2.377469135802469; E=2.807870370370371},
    @{Row=12; B='Claude-3-Haiku'; C=2.361265432098766; D=2.580092592592593; E=2.945061728395062},
    @{Row=13; B='Gemini-2.5-Flash'; C=2.331635802469136; D=2.38858024691358; E=2.631481481481481},
    @{Row=14; B='GPT-5-chat'; C=2.319907407407408; D=2.341975308641976; E=2.731481481481482},
    @{Row=15; B='Claude-3.7-Sonnet'; C=2.312962962962963; D=2.395370370370371; E=2.814043209876543},
    @{Row=16; B='Gemini-2.5-Pro'; C=2.296604938271605; D=2.332407407407407; E=2.890277777777778},
    @{Row=17; B='Mistral-Small-3.1-24b'; C=2.268827160493827; D=1.714814814814815; E=2.637808641975309},
    @{Row=18; B='DeepSeek-Chat-V3.1'; C=2.244753086419753; D=2.315277777777778; E=2.630555555555556},
    @{Row=19; B='Gemini-2.0-Flash-001'; C=2.220833333333333; D=2.223611111111111; E=2.510648148148148},
    @{Row=20; B='Gemini-2.0-Flash-001-Lite'; C=2.217901234567901; D=2.216203703703704; E=2.843827160493827},
    @{Row=21; B='GPT-5-mini'; C=2.206327160493827; D=2.201543209876543; E=2.542746913580247},
    @{Row=22; B='Command-R-08-2024'; C=2.185185185185185; D=2.029938271604939; E=2.481172839506173},
    @{Row=23; B='Gemma-3-27b-It'; C=2.182870370370371; D=2.170833333333333; E=2.737037037037037},
    @{Row=24; B='Llama-3.1-405b'; C=2.181944444444445; D=2.13858024691358; E=2.782253086419753},
    @{Row=25; B='Command-R-Plus-08-2024'; C=2.14429012345679; D=2.029783950617284; E=2.669907407407408},
    @{Row=26; B='GPT-5-nano'; C=2.134876543209876; D=2.165432098765432; E=2.463734567901235},
    @{Row=27; B='GPT-5'; C=2.126234567901235; D=2.173611111111111; E=2.661111111111111},
    @{Row=28; B='Phi-3.5-mini'; C=2.12391975308642; D=2.201851851851852; E=2.577777777777778},
    @{Row=29; B='Llama-3.1-8b'; C=2.118672839506173; D=2.064969135802469; E=2.825154320987654},
    @{Row=30; B='Llama-3.2-90b'; C=2.114660493827161; D=2.090277777777777; E=2.206635802469136},
    @{Row=31; B='Grok-3'; C=2.10570987654321; D=2.139197530864198; E=2.876388888888889},
    @{Row=32; B='Grok-4-Fast'; C=2.093672839506173; D=2.160030864197531; E=2.700308641975309},
    @{Row=33; B='Gemini-2.5-Flash-Lite'; C=2.091203703703703; D=2.175; E=2.633950617283951},
    @{Row=34; B='Llama-3-8b'; C=2.087962962962963; D=1.892901234567901; E=2.745524691358025},
    @{Row=35; B='Claude-3.5-Haiku'; C=2.085493827160493; D=2.400462962962963; E=2.90462962962963},
    @{Row=36; B='Phi-4'; C=2.062345679012346; D=1.758796296296296; E=2.314197530864197},
    @{Row=37; B='Gemma-3-12b-It'; C=2.048765432098766; D=2.912962962962963; E=2.818518518518518},
    @{Row=38; B='Llama-3.3-70b'; C=2.044753086419753; D=2.210802469135802; E=2.724845679012346},
    @{Row=39; B='Qwen3-32b'; C=2.032253086419753; D=1.930246913580247; E=2.430864197530864},
    @{Row=40; B='Llama-4-Maverick'; C=2.024537037037037; D=2.1; E=2.16033950617284},
    @{Row=41; B='Llama-3-70b'; C=2.016975308641975; D=1.99891975308642; E=2.716975308641975},
    @{Row=42; B='Phi-3-mini'; C=2.016666666666667; D=2.150308641975309; E=2.518518518518519},
    @{Row=43; B='Gemma-3n-4B'; C=1.999537037037038; D=2.08425925925926; E=2.604012345679013},
    @{Row=44; B='DeepSeek-Chat-V3-0324'; C=1.983179012345679; D=1.998302469135802; E=2.572530864197531},
    @{Row=45; B='ChatGPT-4o'; C=1.975; D=2.046913580246914; E=2.656481481481482},
    @{Row=46; B='Llama-3.2-3b'; C=1.948302469135802; D=2.043827160493827; E=2.517901234567901},
    @{Row=47; B='Command-A_(Alt)'; C=1.936265432098766; D=1.989506172839506; E=2.415895061728395},
    @{Row=48; B='DeepSeek-V3'; C=1.931172839506172; D=1.973148148148148; E=2.285185185185186},
    @{Row=49; B='GPT-4.1'; C=1.920061728395062; D=1.922376543209876; E=2.571604938271605},
    @{Row=50; B='Command-R7b'; C=1.90679012345679; D=1.971141975308642; E=2.026543209876543},
    @{Row=51; B='Mistral-Medium-3'; C=1.896141975308642; D=1.951388888888889; E=2.672685185185185},
    @{Row=52; B='Phi-3-medium'; C=1.892592592592593; D=1.877777777777778; E=1.833641975308642},
    @{Row=53; B='Llama-4-Scout'; C=1.891820987654321; D=2.020679012345679; E=2.349228395061728},
    @{Row=54; B='Qwen3-8b'; C=1.883796296296297; D=1.947685185185185; E=2.36358024691358},
    @{Row=55; B='Llama-3.2-1b'; C=1.881635802469136; D=1.952314814814815; E=2.579320987654321},
    @{Row=56; B='Qwen3-14b'; C=1.862191358024691; D=1.876234567901234; E=2.363734567901235},
    @{Row=57; B='Gemma-3-4b-It'; C=1.861265432098765; D=2.87962962962963; E=2.819753086419754},
    @{Row=58; B='GPT-4o'; C=1.840895061728395; D=1.84891975308642; E=2.311265432098766},
    @{Row=59; B='Mistral-large-2'; C=1.838117283950617; D=1.84675925925926; E=2.652314814814815},
    @{Row=60; B='GPT-4.1-nano'; C=1.837962962962963; D=1.784104938271605; E=2.043672839506173},
    @{Row=61; B='Gemma-3n-2B'; C=1.821296296296296; D=2.858641975308642; E=2.803703703703704},
    @{Row=62; B='Qwen2.5-72b'; C=1.783796296296296; D=1.781172839506173; E=2.333950617283951},
    @{Row=63; B='Qwen2.5-7b'; C=1.750771604938272; D=1.682870370370371; E=2.039351851851852},
    @{Row=64; B='GPT-4.1-mini'; C=1.719444444444444; D=1.730555555555555; E=2.319135802469136},
    @{Row=65; B='AVG (all models)'; C=2.116174309229865; D=2.139065255731923; E=2.595409563002156}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
